$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4092961847782135
$ws.Range("B1").Value = 1.289472460746765
$ws.Range("C1").Value = 4.512758255004883
$ws.Range("D1").Value = 1.672162413597107
$ws.Range("E1").Value = 1.002083778381348
